$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove six accounts that no longer appear in the export
# (delete from the bottom up so earlier row numbers stay valid):
#   row 13 -> 005110867 DIG            1000
#   row 12 -> 004224815 GUILHERME      1500
#   row 7  -> 004751770 DILSON         15315.3
#   row 6  -> 004200433 BENTO          15409.32
#   row 4  -> 004890544 ASSAKO         18970.25
#   row 3  -> 004335144 EDMUNDO        20735.88
$ws.Rows(13).Delete()
$ws.Rows(12).Delete()
$ws.Rows(7).Delete()
$ws.Rows(6).Delete()
$ws.Rows(4).Delete()
$ws.Rows(3).Delete()

# Insert the new account row right after ANTONIO (which is now row 23,
# having shifted up from row 29 after the six deletions above) so it
# lands before NATHALIA.
$ws.Rows(24).Insert()

# Keep the leading zeros in the account number by forcing the cell to
# Text before writing the value (otherwise Excel would coerce the
# numeric-looking string to a number and drop them).
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "005611045"
$ws.Range("B24").Value = "PAULA"
$ws.Range("C24").Value = 100
